# Saldo_guide.xlsx update
# - Roll the "G" (date) column forward one day for every data row (45418 -> 45419)
# - Re-settle a handful of rows where the pending ("D") amount was fully applied,
#   so the remainder ("E") drops to 0 and D absorbs the total (D+E, i.e. column H)
# - One row (169) goes the other way: a new pending amount lands in E, so H grows
# - Leave the active sheet fully selected (as if the user pressed Ctrl+A last)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: shift every data row's date from 45418 to 45419 -------------
for ($r = 2; $r -le 311; $r++) {
    $ws.Cells.Item($r, 7).Value = 45419
}

# --- Column D/E (and, where it actually changes, H) rebalancing ------------
$ws.Cells.Item(2, 4).Value = 1972.16
$ws.Cells.Item(2, 5).Value = 0

$ws.Cells.Item(5, 4).Value = 43402.17
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 8).Value = 43402.17

$ws.Cells.Item(12, 4).Value = 4510.54
$ws.Cells.Item(12, 5).Value = 0

$ws.Cells.Item(23, 4).Value = 12659.96
$ws.Cells.Item(23, 5).Value = 0

$ws.Cells.Item(27, 4).Value = 4298.17
$ws.Cells.Item(27, 5).Value = 0

$ws.Cells.Item(55, 4).Value = 5369.68
$ws.Cells.Item(55, 5).Value = 0

$ws.Cells.Item(61, 4).Value = 63487.69
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 8).Value = 63487.69

$ws.Cells.Item(71, 4).Value = 1761.17
$ws.Cells.Item(71, 5).Value = 0

$ws.Cells.Item(73, 4).Value = 5248.07
$ws.Cells.Item(73, 5).Value = 0

$ws.Cells.Item(115, 4).Value = 5844.37
$ws.Cells.Item(115, 5).Value = 0

$ws.Cells.Item(120, 4).Value = 15901.01
$ws.Cells.Item(120, 5).Value = 0

$ws.Cells.Item(125, 4).Value = 18770.07
$ws.Cells.Item(125, 5).Value = 0

$ws.Cells.Item(151, 4).Value = 2984.41
$ws.Cells.Item(151, 5).Value = 0

$ws.Cells.Item(167, 4).Value = 72750.07
$ws.Cells.Item(167, 5).Value = 0

$ws.Cells.Item(169, 5).Value = 2083.6999999999998
$ws.Cells.Item(169, 8).Value = 2134.15

$ws.Cells.Item(184, 4).Value = 1064.5
$ws.Cells.Item(184, 5).Value = 0

$ws.Cells.Item(200, 4).Value = 1257.18
$ws.Cells.Item(200, 5).Value = 0

$ws.Cells.Item(267, 4).Value = 2013.32
$ws.Cells.Item(267, 5).Value = 0

$ws.Cells.Item(274, 4).Value = 5502.35
$ws.Cells.Item(274, 5).Value = 0

$ws.Cells.Item(290, 4).Value = 6123.43
$ws.Cells.Item(290, 5).Value = 0

# --- Selection: select the whole sheet (Ctrl+A), matching the saved state --
$ws.Cells.Select()
